$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived NATMI values for the Efna3-Epha2 sheet, plus the new
# "MuSCs" sending-cluster rows (12-16) that the refreshed script now emits.
$data = @(
    @("ECs", "Efna3", "Epha2", "ECs", 2, 0.6666666666666666, 0.2420556666666667, 0.726167, 0.5314769098578004, 0.5314769098578004, 3, 1, 20.66830833333333, 62.004925, 0.6755285375771634, 0.6755285375771634, 5.002881152497222, 45.025930372475, 0.3590278196722698, 0.3590278196722698),
    @("ECs", "Efna3", "Epha2", "FAPs", 2, 0.6666666666666666, 0.2420556666666667, 0.726167, 0.5314769098578004, 0.5314769098578004, 3, 1, 2.156330666666667, 6.468992, 0.07047809033489469, 0.07047809033489467, 0.5219520570737778, 4.697568513664, 0.03745747766386873, 0.03745747766386873),
    @("ECs", "Efna3", "Epha2", "Inflammatory-Mac", 2, 0.6666666666666666, 0.2420556666666667, 0.726167, 0.5314769098578004, 0.5314769098578004, 3, 1, 1.683564, 5.050692, 0.05502605769642779, 0.05502605769642779, 0.4075162063959999, 3.667645857564, 0.02924507910615448, 0.02924507910615448),
    @("ECs", "Efna3", "Epha2", "MuSCs", 2, 0.6666666666666666, 0.2420556666666667, 0.726167, 0.5314769098578004, 0.5314769098578004, 3, 1, 5.278649666666666, 15.835949, 0.1725288026574751, 0.1725288026574751, 1.277727064164778, 11.499543577483, 0.09169507489786112, 0.0916950748978611),
    @("ECs", "Efna3", "Epha2", "Resolving-Mac", 2, 0.6666666666666666, 0.2420556666666667, 0.726167, 0.5314769098578004, 0.5314769098578004, 3, 1, 0.8089063333333334, 2.426719, 0.02643851173403914, 0.02643851173403913, 0.1958003617858889, 1.762203256073, 0.01405145851764632, 0.01405145851764631),
    @("FAPs", "Efna3", "Epha2", "ECs", 1, 0.3333333333333333, 0.2054156666666667, 0.616247, 0.4510271759376837, 0.4510271759376837, 3, 1, 20.66830833333333, 62.004925, 0.6755285375771634, 0.6755285375771634, 4.245594335163888, 38.210349016475, 0.3046817285687414, 0.3046817285687414),
    @("FAPs", "Efna3", "Epha2", "FAPs", 1, 0.3333333333333333, 0.2054156666666667, 0.616247, 0.4510271759376837, 0.4510271759376837, 3, 1, 2.156330666666667, 6.468992, 0.07047809033489469, 0.07047809033489467, 0.4429441014471112, 3.986496913024, 0.03178753404922851, 0.0317875340492285),
    @("FAPs", "Efna3", "Epha2", "Inflammatory-Mac", 1, 0.3333333333333333, 0.2054156666666667, 0.616247, 0.4510271759376837, 0.4510271759376837, 3, 1, 1.683564, 5.050692, 0.05502605769642779, 0.05502605769642779, 0.345830421436, 3.112473792924, 0.02481824740580387, 0.02481824740580387),
    @("FAPs", "Efna3", "Epha2", "MuSCs", 1, 0.3333333333333333, 0.2054156666666667, 0.616247, 0.4510271759376837, 0.4510271759376837, 3, 1, 5.278649666666666, 15.835949, 0.1725288026574751, 0.1725288026574751, 1.084317340378111, 9.758856063403, 0.07781517863051092, 0.0778151786305109),
    @("FAPs", "Efna3", "Epha2", "Resolving-Mac", 1, 0.3333333333333333, 0.2054156666666667, 0.616247, 0.4510271759376837, 0.4510271759376837, 3, 1, 0.8089063333333334, 2.426719, 0.02643851173403914, 0.02643851173403913, 0.1661620337325556, 1.495458303593, 0.01192448728339898, 0.01192448728339898),
    @("MuSCs", "Efna3", "Epha2", "ECs", 1, 0.3333333333333333, 0.007968333333333332, 0.023905, 0.01749591420451593, 0.01749591420451593, 3, 1, 20.66830833333333, 62.004925, 0.6755285375771634, 0.6755285375771634, 0.1646919702361111, 1.482227732125, 0.01181898933615216, 0.01181898933615216),
    @("MuSCs", "Efna3", "Epha2", "FAPs", 1, 0.3333333333333333, 0.007968333333333332, 0.023905, 0.01749591420451593, 0.01749591420451593, 3, 1, 2.156330666666667, 6.468992, 0.07047809033489469, 0.07047809033489467, 0.01718236152888889, 0.15464125376, 0.001233078621797441, 0.00123307862179744),
    @("MuSCs", "Efna3", "Epha2", "Inflammatory-Mac", 1, 0.3333333333333333, 0.007968333333333332, 0.023905, 0.01749591420451593, 0.01749591420451593, 3, 1, 1.683564, 5.050692, 0.05502605769642779, 0.05502605769642779, 0.01341519914, 0.12073679226, 0.0009627311844694439, 0.0009627311844694439),
    @("MuSCs", "Efna3", "Epha2", "MuSCs", 1, 0.3333333333333333, 0.007968333333333332, 0.023905, 0.01749591420451593, 0.01749591420451593, 3, 1, 5.278649666666666, 15.835949, 0.1725288026574751, 0.1725288026574751, 0.04206204009388888, 0.378558360845, 0.003018549129103044, 0.003018549129103043),
    @("MuSCs", "Efna3", "Epha2", "Resolving-Mac", 1, 0.3333333333333333, 0.007968333333333332, 0.023905, 0.01749591420451593, 0.01749591420451593, 3, 1, 0.8089063333333334, 2.426719, 0.02643851173403914, 0.02643851173403913, 0.006445635299444444, 0.05801071769500001, 0.0004625659329938363, 0.0004625659329938363),
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}